$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44159
$ws.Range("M2").Value = 320
$ws.Range("N2").Value = 13500
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 13750
$ws.Range("S2").Value = 1719
# Row 3
$ws.Range("D3").Value = 44159
$ws.Range("N3").Value = 11500
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 11750
$ws.Range("S3").Value = 1469
# Row 4
$ws.Range("D4").Value = 44169
$ws.Range("M4").Value = 240
$ws.Range("N4").Value = 14500
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14750
$ws.Range("S4").Value = 1844
# Row 5
$ws.Range("D5").Value = 44169
$ws.Range("M5").Value = 240
$ws.Range("N5").Value = 12500
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 12750
$ws.Range("S5").Value = 1594
# Row 6
$ws.Range("D6").Value = 44169
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10500
$ws.Range("P6").Value = 10250
$ws.Range("S6").Value = 1281
# Row 7
$ws.Range("D7").Value = 44449
$ws.Range("M7").Value = 240
$ws.Range("N7").Value = 2900
$ws.Range("O7").Value = 3000
$ws.Range("P7").Value = 2950
$ws.Range("Q7").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("S7").Value = 2950
$ws.Range("T7").Value = 1
# Row 8
$ws.Range("D8").Value = 44449
$ws.Range("L8").Value = "Extra (doble especial)"
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 3100
$ws.Range("O8").Value = 3200
$ws.Range("P8").Value = 3150
$ws.Range("Q8").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("S8").Value = 3150
$ws.Range("T8").Value = 1
# Row 9
$ws.Range("D9").Value = 44449
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 2700
$ws.Range("O9").Value = 2800
$ws.Range("P9").Value = 2750
$ws.Range("Q9").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("S9").Value = 2750
$ws.Range("T9").Value = 1
# Row 10
$ws.Range("D10").Value = 44161
$ws.Range("N10").Value = 13500
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 13750
$ws.Range("S10").Value = 1719
# Row 11
$ws.Range("D11").Value = 44161
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 11500
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 11750
$ws.Range("S11").Value = 1469
# Row 12
$ws.Range("D12").Value = 44161
$ws.Range("N12").Value = 9000
$ws.Range("O12").Value = 9500
$ws.Range("P12").Value = 9250
$ws.Range("S12").Value = 1156
# Row 13
$ws.Range("D13").Value = 44165
$ws.Range("M13").Value = 300
# Row 14
$ws.Range("D14").Value = 44165
$ws.Range("M14").Value = 240
# Row 15
$ws.Range("D15").Value = 44165
# Row 16
$ws.Range("D16").Value = 44172
$ws.Range("L16").Value = "Especial"
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 14500
$ws.Range("P16").Value = 14250
$ws.Range("Q16").Value = "`$/bandeja 8 kilos"
$ws.Range("R16").Value = "Provincia de Limarí"
$ws.Range("S16").Value = 1781
$ws.Range("T16").Value = 8
# Row 17
$ws.Range("D17").Value = 44172
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12500
$ws.Range("P17").Value = 12250
$ws.Range("Q17").Value = "`$/bandeja 8 kilos"
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("S17").Value = 1531
$ws.Range("T17").Value = 8
# Row 18
$ws.Range("D18").Value = 44172
$ws.Range("L18").Value = "Segunda"
$ws.Range("N18").Value = 9500
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 9750
$ws.Range("Q18").Value = "`$/bandeja 8 kilos"
$ws.Range("R18").Value = "Provincia de Limarí"
$ws.Range("S18").Value = 1219
$ws.Range("T18").Value = 8
# Row 19
$ws.Range("D19").Value = 44166
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 14000
$ws.Range("O19").Value = 14500
$ws.Range("P19").Value = 14250
$ws.Range("S19").Value = 1781
# Row 20
$ws.Range("D20").Value = 44166
$ws.Range("N20").Value = 12000
$ws.Range("O20").Value = 12500
$ws.Range("P20").Value = 12250
$ws.Range("S20").Value = 1531
# Row 24
$ws.Range("D24").Value = 44162
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 340
$ws.Range("N24").Value = 14000
$ws.Range("O24").Value = 14500
$ws.Range("P24").Value = 14250
$ws.Range("Q24").Value = "`$/bandeja 8 kilos"
$ws.Range("R24").Value = "Provincia de Limarí"
$ws.Range("S24").Value = 1781
$ws.Range("T24").Value = 8
# Row 25
$ws.Range("D25").Value = 44162
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 300
$ws.Range("N25").Value = 12000
$ws.Range("O25").Value = 12500
$ws.Range("P25").Value = 12250
$ws.Range("Q25").Value = "`$/bandeja 8 kilos"
$ws.Range("R25").Value = "Provincia de Limarí"
$ws.Range("S25").Value = 1531
$ws.Range("T25").Value = 8
# Row 26
$ws.Range("D26").Value = 44162
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 200
$ws.Range("N26").Value = 9500
$ws.Range("O26").Value = 10000
$ws.Range("P26").Value = 9750
$ws.Range("Q26").Value = "`$/bandeja 8 kilos"
$ws.Range("S26").Value = 1219
$ws.Range("T26").Value = 8
# Row 27
$ws.Range("D27").Value = 44410
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 240
$ws.Range("N27").Value = 2400
$ws.Range("O27").Value = 2500
$ws.Range("P27").Value = 2450
$ws.Range("R27").Value = "Provincia del Elquí"
$ws.Range("S27").Value = 2450
# Row 28
$ws.Range("D28").Value = 44410
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 240
$ws.Range("N28").Value = 2000
$ws.Range("O28").Value = 2100
$ws.Range("P28").Value = 2050
$ws.Range("R28").Value = "Provincia del Elquí"
$ws.Range("S28").Value = 2050
# Row 29
$ws.Range("D29").Value = 44410
$ws.Range("L29").Value = "Tercera"
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = 1600
$ws.Range("O29").Value = 1700
$ws.Range("P29").Value = 1650
$ws.Range("Q29").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R29").Value = "Provincia del Elquí"
$ws.Range("S29").Value = 1650
$ws.Range("T29").Value = 1
# Row 30
$ws.Range("D30").Value = 44411
$ws.Range("M30").Value = 600
$ws.Range("N30").Value = 2400
$ws.Range("O30").Value = 2500
$ws.Range("P30").Value = 2450
$ws.Range("Q30").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R30").Value = "Provincia del Elquí"
$ws.Range("S30").Value = 2450
$ws.Range("T30").Value = 1
# Row 31
$ws.Range("D31").Value = 44411
$ws.Range("M31").Value = 400
$ws.Range("N31").Value = 2000
$ws.Range("O31").Value = 2100
$ws.Range("P31").Value = 2050
$ws.Range("Q31").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R31").Value = "Provincia del Elquí"
$ws.Range("S31").Value = 2050
$ws.Range("T31").Value = 1
